$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 1.1.0 -> 1.2.0
$ws.Range("B3").Value = "1.2.0"

# Title value (B5) was empty -> same as Name (B4) = "LocallyDefinedInterventions"
$ws.Range("B5").Value = "LocallyDefinedInterventions"

# Experimental value (B7) was empty -> "false"
$ws.Range("B7").Value = "'false"

# Date: 2023-06-21T21:59:46+02:00 -> 2024-10-31T19:21:51+01:00
$ws.Range("B8").Value = "2024-10-31T19:21:51+01:00"

# Contact value (B10): "No display for ContactDetail" -> "KL (http://www.kl.dk)"
$ws.Range("B10").Value = "KL (http://www.kl.dk)"

# Jurisdiction value (B11): "Denmark" -> "" (empty)
$ws.Range("B11").Value = ""

# Case Sensitive value (B15) was empty -> "true"
$ws.Range("B15").Value = "'true"
